$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old email everywhere it is used (B2, B3) with the new one
$ws.Range("B2").Value = "biktoras92@gmail.com"
$ws.Range("B3").Value = "biktoras92@gmail.com"

# B4 gets a brand new distinct email address
$ws.Range("B4").Value = "victorsfak03@gmail.com"

# Row 5 gets filled in with a duplicate "BC" name and the recurring email
$ws.Range("A5").Value = "BC ΧΑΤΖΗΔΑΚΙ ΣΤΥΛΙΑΝΗ"
$ws.Range("B5").Value = "biktoras92@gmail.com"
